$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overal Stats")
$ws.Range("A1").Value = "test"
